# Fruta / hortaliza, semanal
# Insert two new weekly data rows (date 44610) right above the existing
# block of rows, shifting all rows from 542 onward down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 542:543 (pushes old row 542.. down to 544..)
$ws.Rows("542:543").Insert()

# Row 542 - "Primera" quality record for 2022-02-18 (serial 44610)
$ws.Cells.Item(542, 1).Value = 8
$ws.Cells.Item(542, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(542, 3).Value = "Coquimbo"
$ws.Cells.Item(542, 4).Value = 44610
$ws.Cells.Item(542, 5).Value = 4
$ws.Cells.Item(542, 6).Value = 100112008
$ws.Cells.Item(542, 7).Value = "Coliflor"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 2400
$ws.Cells.Item(542, 11).Value = 850
$ws.Cells.Item(542, 12).Value = 900
$ws.Cells.Item(542, 13).Value = 875
$ws.Cells.Item(542, 14).Value = "$/unidad"
$ws.Cells.Item(542, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(542, 16).Value = 875
$ws.Cells.Item(542, 17).Value = 1
$ws.Cells.Item(542, 18).Value = "Hortaliza"

# Row 543 - "Segunda" quality record for 2022-02-18 (serial 44610)
$ws.Cells.Item(543, 1).Value = 8
$ws.Cells.Item(543, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(543, 3).Value = "Coquimbo"
$ws.Cells.Item(543, 4).Value = 44610
$ws.Cells.Item(543, 5).Value = 4
$ws.Cells.Item(543, 6).Value = 100112008
$ws.Cells.Item(543, 7).Value = "Coliflor"
$ws.Cells.Item(543, 8).Value = "Sin especificar"
$ws.Cells.Item(543, 9).Value = "Segunda"
$ws.Cells.Item(543, 10).Value = 1300
$ws.Cells.Item(543, 11).Value = 750
$ws.Cells.Item(543, 12).Value = 800
$ws.Cells.Item(543, 13).Value = 775
$ws.Cells.Item(543, 14).Value = "$/unidad"
$ws.Cells.Item(543, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(543, 16).Value = 775
$ws.Cells.Item(543, 17).Value = 1
$ws.Cells.Item(543, 18).Value = "Hortaliza"
